# Update "想去人数" (F) and "最低票价" (G) values on the "展览" and "全部类型"
# sheets to reflect the latest scrape output.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of cell address -> new value
$updates = @{
    "G3"  = 70
    "F6"  = 1350
    "F7"  = 1582
    "F8"  = 348
    "F9"  = 451
    "F11" = 182
    "F15" = 287
    "F17" = 333
    "F21" = 182
    "F22" = 700
    "F25" = 4291
    "F27" = 295
    "F28" = 1132
    "F29" = 497
    "F31" = 647
    "F33" = 328
    "F35" = 168
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
